$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.675.64'
$ws.Range('E2').Value = '  +1.28%  '
$ws.Range('D3').Value = '3.579.25'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('E4').Value = '  +0.00%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '589.25'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +2.48%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '187.01'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +0.94%  '
$ws.Range('D7').Value = '3.568.82'
$ws.Range('E7').Value = '  +0.43%  '
$ws.Range('E8').Value = '  +0.96%  '
$ws.Range('E9').Value = '  +0.10%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.202'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +10.40%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.650'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +0.63%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '54.51'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -0.14%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '0.0000313'
$c.Style = "Normal"
$ws.Range('E13').Value = '  +4.19%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '9.56'
$c.Style = "Normal"
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('D15').Value = '4.152.46'
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('D17').Value = '70.676.02'
$ws.Range('E17').Value = '  +1.33%  '
$ws.Range('D18').Value = '3.567.45'
$ws.Range('E18').Value = '  +0.31%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '12.47'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('E20').Value = '  +0.05%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '561.95'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +13.80%  '
$ws.Range('E22').Value = '  -0.68%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '17.84'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -6.25%  '
$ws.Range('E24').Value = '  +7.70%  '
$ws.Range('E25').Value = '  +0.70%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '95.75'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +0.69%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '11.54'
$c.Style = "Normal"
$ws.Range('E27').Value = '  +1.75%  '
$ws.Range('E28').Value = '  +1.70%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '9.16'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -1.02%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '32.22'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +2.21%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '7.27'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -3.01%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '12.48'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +4.39%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '65.16'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -2.37%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.116'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +1.07%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '565.12'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +0.07%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '3.31'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +6.10%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.418'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +7.04%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '38.10'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -0.95%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('D40').Value = '0.0₃0777'
$ws.Range('E40').Value = '  -0.94%  '
$ws.Range('E41').Value = '  +1.14%  '
$ws.Range('D42').Value = '3.345.87'
$ws.Range('E42').Value = '  +3.92%  '
$ws.Range('E43').Value = '  -3.90%  '
$ws.Range('E44').Value = '  -2.22%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '3.59'
$c.Style = "Normal"
$ws.Range('E45').Value = '  +3.96%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.0445'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +1.64%  '
$ws.Range('B47').Value = 'ThetaToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '2.97'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('E48').Value = '  -1.21%  '
$ws.Range('E49').Value = '  +1.40%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -0.11%  '
$ws.Range('E51').Value = '  +20.51%  '
